$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds date-serial values. Every populated data row
# (C2:C103) currently stores 45177 (2023-09-08) and must be bumped by one
# day to 45178 (2023-09-09).
for ($row = 2; $row -le 103; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45177) {
        $cell.Value2 = 45178
    }
}
